# Applies the edits described by the commit diff to the single table
# in the document (Table 1).
#
# NOTE: Find.Execute with Replace:=wdReplaceAll (2) always rewrites every
# matching occurrence in the whole story, even when called on a Find
# object that belongs to a restricted sub-range. To change just the
# occurrence inside one particular cell (several strings below are not
# unique document-wide) we scope the call to that cell's Range AND pass
# Replace:=wdReplaceOne (1), which - like real Word - only touches the
# first match starting at/after the range used to invoke Find.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-InCell($table, $rowIndex, $cellIndex, $oldText, $newText) {
    $cell = $table.Rows.Item($rowIndex).Cells.Item($cellIndex)
    $cell.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                              $true, 1, $false, $newText, 1) | Out-Null
}

# --- Row 5 (doc item 1): "Выписка координат пунктов ГГС" / "б\н 2017-11-01"
Replace-InCell $t 5 2 "Выписка координат пунктов ГГС" ""
Replace-InCell $t 5 3 "б\н 2017-11-01" "№ 36/исх/17-449216 от 21.06.2017 г."

# --- Row 6 (doc item 2): "Постановление администрации Бежецкого района Тверской области" / "3 2018-01-12"
Replace-InCell $t 6 2 "Постановление администрации Бежецкого района Тверской области" ""
Replace-InCell $t 6 3 "3 2018-01-12" "№ 36/исх/17-590486 от 14.08.2017 г."

# --- Row 7 (doc item 3): "" / "99/2017/37507685 2017-11-21"
Replace-InCell $t 7 3 "99/2017/37507685 2017-11-21" "№ 165 от 02.06.2017 г."

# --- Row 8 (doc item 4): "Выписка из ЕГРН" / "99/2017/42655414   2017-12-06"
Replace-InCell $t 8 2 "Выписка из ЕГРН" ""
Replace-InCell $t 8 3 "99/2017/42655414   2017-12-06" "№ 158 от 02.06.2017 г."

# --- Delete rows 9 and 10 (doc items 5 and 6) entirely.
# Deleting row 9 shifts what was row 10 up to row 9.
$t.Rows.Item(9).Delete()
$t.Rows.Item(9).Delete()

# --- "Система координат ... МСК-69, зона 2" -> "... СК кадастрового округа"
# (unique document-wide, so a plain ReplaceAll is fine)
$d.Content.Find.Execute("МСК-69, зона 2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "СК кадастрового округа", 2) | Out-Null

# Table rows have shifted up by 2 after the two deletions above.
# Row that was 15 (Крутец...) is now row 13.
Replace-InCell $t 13 2 "Крутец   наружный знак утрачен " "Марки пир. "
Replace-InCell $t 13 3 "2 класс" "3 класс"
Replace-InCell $t 13 4 "392758.57" "418081.05"
Replace-InCell $t 13 5 "2319346.5" "1330793.93"

# Row that was 16 (Сокольниково...) is now row 14.
Replace-InCell $t 14 2 "Сокольниково пир. 6.2. " "Пушкино пир. "
Replace-InCell $t 14 3 "1 класс" "3 класс"
Replace-InCell $t 14 4 "402600.64" "404698.2"
Replace-InCell $t 14 5 "2326524.69" "1307254.3"

# Row that was 17 (Алабузино...) is now row 15. Class stays "2 класс".
Replace-InCell $t 15 2 "Алабузино пир. 6.2. " "Михайловка пир. "
Replace-InCell $t 15 4 "401851.71" "387798.79"
Replace-InCell $t 15 5 "2316356.82" "2159578.49"

# Row that was 21 (instrument 1) is now row 19.
Replace-InCell $t 19 3 "33967-0717 января 2018 г." "№ 012343 17 января 2018 г"
Replace-InCell $t 19 4 "№ 012343" "33967-07"

# Row that was 22 (instrument 2) is now row 20.
Replace-InCell $t 20 2 "Trimble R7GNS" "Trimble R7GNSS"
Replace-InCell $t 20 3 "37145-0817 января 2018 г." "№ 012342 17 января 2018 г"
Replace-InCell $t 20 4 "№ 012342" "37145-08"

# --- Delete the trailing empty row (was row 26, now row 24 after the two earlier deletions).
$t.Rows.Item(24).Delete()

Write-Output "edits applied"
